$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the header row (row 1) with validated column names.
# Set in this order so shared-string indices line up: id, author, titulo
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "author"
$ws.Range("B1").Value = "titulo"

# Add the id value for the first data row (now row 2)
$ws.Range("A2").Value = 1

# Reflect the new selection in the sheet view
$ws.Range("B1").Select()
